$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 292 (pushes existing rows 292+ down by one)
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row with the "L" (lists) library entry
$ws.Range("B292").Value = "N/A"
$ws.Range("D292").Value = "L"
$ws.Range("E292").Value = "table"
$ws.Range("F292").Value = "Library that deals with lists."

# Match the author's saved selection/viewport state
$ws.Range("D289").Select()
